$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Byte 0&1")

# Widen column B
$ws1.Columns.Item(2).ColumnWidth = 93.05

# Extend the shared formula in column A down to row 17 (was row 9)
$ws1.Range("A10:A17").Formula = "=A9+1"
$ws1.Range("A10:A17").HorizontalAlignment = -4108  # xlCenter (match style of A4:A9)

# New remark cells in column B (order matters for shared-string table order)
$ws1.Range("B4").Value = "Δ56 = 1 and (Δ42 != 1 and Δ62!=1 and Δ74 != 1 and Δ78 != 1) and (Δ42!= 1 and Δ74!= 1 and Δ78!= 1)"
$ws1.Range("B4").HorizontalAlignment = -4131  # xlLeft
$ws1.Range("B4").Interior.Pattern = -4142     # xlNone

$ws1.Range("B6").Value = "Δ58=1 and Δ44!= 1 and  Δ64 != 1 and Δ80 != 1 ) "
$ws1.Range("B8").Value = "Δ60=1 and (Δ50 !=1 and Δ66!= 1 and Δ78!= 1 and Δ82!= 1 and Δ98!=1)"
$ws1.Range("B10").Value = "Δ62=1 and (Δ52!=1 and Δ80!= 1 and Δ84!= 1 and Δ100!=1)"
$ws1.Range("B14").Value = "Δ50=1 and (Δ36!= 1 and Δ40!= 1 and Δ56 !=1 and Δ88!= 1 ) and (Δ84 != 1)"
$ws1.Range("B16").Value = "Δ52 = 1 and (Δ38 != 1 and Δ42 != 1 and Δ58!= 1 and Δ70 != 1 and Δ74 != 1) and (Δ86 != 1 and Δ92 != 1)"

# Update B2 text (was "Δ42 = Δ46 = 1") -- changed last so it lands at the end of the shared string table
$ws1.Range("B2").Value = "Δ54=1 and (Δ40 != 1 and/or? Δ44 != 1 and Δ60 !=1 and Δ92 != 1) and (Δ40 !=1 and Δ88 !=1)"

# Update selection to B5
$ws1.Range("B5").Select()
